# Natmi following Dr Hou advice
# Rebuild the LR-pair table for Sfrp1-Fzd2 across 3 sending/target clusters (ECs, FAPs, sCs)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "sCs")

# Per-cluster ligand stats (Sfrp1), keyed by cluster: E,F,G,H,I,J
$ligandStats = @{
    "ECs"  = @(1, 0.3333333333333333, 0.02759233333333333, 0.082777, 0.003002532875897786, 0.003002532875897787)
    "FAPs" = @(3, 1, 7.96874, 23.90622, 0.8671395615744129, 0.867139561574413)
    "sCs"  = @(3, 1, 1.193353333333333, 3.58006, 0.1298579055496893, 0.1298579055496893)
}

# Per-cluster receptor stats (Fzd2), keyed by cluster: K,L,M,N,O,P
$receptorStats = @{
    "ECs"  = @(1, 0.3333333333333333, 0.1278803333333333, 0.383641, 0.009974564977605908, 0.009974564977605908)
    "FAPs" = @(3, 1, 10.61985133333333, 31.859554, 0.8283400145723324, 0.8283400145723324)
    "sCs"  = @(3, 1, 2.072911, 6.218733, 0.1616854204500617, 0.1616854204500617)
}

# Per sending/target cluster pair edge stats: Q,R,S,T
$edgeStats = @{
    "ECs|ECs"   = @(0.003528516784111111, 0.031756651057, 0.0000299489592680404, 0.00002994895926804041)
    "ECs|FAPs"  = @(0.2930264779397778, 2.637238301458, 0.002487118126175079, 0.00248711812617508)
    "ECs|sCs"   = @(0.05719645128233333, 0.5147680615410001, 0.0004854657904546664, 0.0004854657904546665)
    "FAPs|ECs"  = @(1.019045127446667, 9.171406147020001, 0.00864933990157668, 0.00864933990157668)
    "FAPs|FAPs" = @(84.62683411398667, 761.6415070258801, 0.7182863970707951, 0.7182863970707952)
    "FAPs|sCs"  = @(16.51848880214, 148.66639921926, 0.1402038246020411, 0.1402038246020411)
    "sCs|ECs"   = @(0.1526064220511111, 1.37345779846, 0.001295276116761186, 0.001295276116761186)
    "sCs|FAPs"  = @(12.67323498813778, 114.05911489324, 0.1075664993753622, 0.1075664993753622)
    "sCs|sCs"   = @(2.473715251553333, 22.26343726398, 0.02099613005756591, 0.02099613005756591)
}

$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ls = $ligandStats[$sending]
        $rs = $receptorStats[$target]
        $es = $edgeStats["$sending|$target"]

        $ws.Cells.Item($row, 1).Value = $sending
        $ws.Cells.Item($row, 2).Value = "Sfrp1"
        $ws.Cells.Item($row, 3).Value = "Fzd2"
        $ws.Cells.Item($row, 4).Value = $target

        $ws.Cells.Item($row, 5).Value = $ls[0]
        $ws.Cells.Item($row, 6).Value = $ls[1]
        $ws.Cells.Item($row, 7).Value = $ls[2]
        $ws.Cells.Item($row, 8).Value = $ls[3]
        $ws.Cells.Item($row, 9).Value = $ls[4]
        $ws.Cells.Item($row, 10).Value = $ls[5]

        $ws.Cells.Item($row, 11).Value = $rs[0]
        $ws.Cells.Item($row, 12).Value = $rs[1]
        $ws.Cells.Item($row, 13).Value = $rs[2]
        $ws.Cells.Item($row, 14).Value = $rs[3]
        $ws.Cells.Item($row, 15).Value = $rs[4]
        $ws.Cells.Item($row, 16).Value = $rs[5]

        $ws.Cells.Item($row, 17).Value = $es[0]
        $ws.Cells.Item($row, 18).Value = $es[1]
        $ws.Cells.Item($row, 19).Value = $es[2]
        $ws.Cells.Item($row, 20).Value = $es[3]

        $row = $row + 1
    }
}
